$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036217019838033
$bf[0,2] = 1.039418025472014
$bf[0,3] = 1.049455627872843
$bf[0,4] = 1.056295776982071
$ws.Range("B2:F2").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037466682437606
$in[0,1] = 1.041326542649597
$in[0,2] = 1.042203270430258
$in[0,3] = 1.052212637466212
$in[0,4] = 1.059033878659115
$in[0,5] = 1.017709500113599
$ws.Range("I2:N2").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037085915620453
$bf[0,2] = 1.040069383033446
$bf[0,3] = 1.050344663540063
$bf[0,4] = 1.057326228297207
$ws.Range("B3:F3").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037655478765406
$in[0,1] = 1.041839716938501
$in[0,2] = 1.042665261404103
$in[0,3] = 1.05291372443103
$in[0,4] = 1.059877388718409
$in[0,5] = 1.017881637530269
$ws.Range("I3:N3").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037648672763576
$bf[0,2] = 1.040491227766665
$bf[0,3] = 1.050920868549684
$bf[0,4] = 1.057994263567493
$ws.Range("B4:F4").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.03777668340538
$in[0,1] = 1.042171646524717
$in[0,2] = 1.04296388741798
$in[0,3] = 1.053367676425191
$in[0,4] = 1.060423850578005
$in[0,5] = 1.017992928876836
$ws.Range("I4:N4").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037885379886713
$bf[0,2] = 1.04066865901973
$bf[0,3] = 1.051163328559539
$bf[0,4] = 1.058275406843219
$ws.Range("B5:F5").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037827407612378
$in[0,1] = 1.042311157919192
$in[0,2] = 1.043089353776195
$in[0,3] = 1.053558588861722
$in[0,4] = 1.060653738498734
$in[0,5] = 1.018039693073232
$ws.Range("I5:N5").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037925131289159
$bf[0,2] = 1.040698455607707
$bf[0,3] = 1.051204051720465
$bf[0,4] = 1.058322629659773
$ws.Range("B6:F6").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037835910912477
$in[0,1] = 1.04233458060097
$in[0,2] = 1.043110415640825
$in[0,3] = 1.053590648028277
$in[0,4] = 1.060692346763942
$in[0,5] = 1.018047543645537
$ws.Range("I6:N6").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037651835170044
$bf[0,2] = 1.040493598268967
$bf[0,3] = 1.050924107436163
$bf[0,4] = 1.057998019035803
$ws.Range("B7:F7").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037777362090067
$in[0,1] = 1.042173510809066
$in[0,2] = 1.042965564205644
$in[0,3] = 1.053370227129013
$in[0,4] = 1.060426921742863
$in[0,5] = 1.017993553831962
$ws.Range("I7:N7").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036510558065185
$bf[0,2] = 1.039638076691508
$bf[0,3] = 1.049755886215468
$bf[0,4] = 1.056643760164238
$ws.Range("B8:F8").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037530685297382
$in[0,1] = 1.041499998281799
$in[0,2] = 1.042359466612472
$in[0,3] = 1.052449509839828
$in[0,4] = 1.059318810529987
$in[0,5] = 1.017767693787905
$ws.Range("I8:N8").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034503548674164
$bf[0,2] = 1.038133461459059
$bf[0,3] = 1.047704592443709
$bf[0,4] = 1.054267130433082
$ws.Range("B9:F9").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037088688891084
$in[0,1] = 1.040312245124468
$in[0,2] = 1.041289095281954
$in[0,3] = 1.050829455197794
$in[0,4] = 1.057371247747656
$in[0,5] = 1.017369004555578
$ws.Range("I9:N9").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033168353533969
$bf[0,2] = 1.037132440965156
$bf[0,3] = 1.046342033799903
$bf[0,4] = 1.052689349736648
$ws.Range("B10:F10").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036789137746173
$in[0,1] = 1.039519840815798
$in[0,2] = 1.040573998334805
$in[0,3] = 1.049751083940311
$in[0,4] = 1.056076360144942
$in[0,5] = 1.017102768435459
$ws.Range("I10:N10").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032590882061627
$bf[0,2] = 1.036699494919174
$bf[0,3] = 1.045753228322652
$bf[0,4] = 1.052007745808406
$ws.Range("B11:F11").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036658276999531
$in[0,1] = 1.039176599000597
$in[0,2] = 1.040264007265977
$in[0,3] = 1.049284547667994
$in[0,4] = 1.055516503061359
$in[0,5] = 1.016987385132502
$ws.Range("I11:N11").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032376486392703
$bf[0,2] = 1.036538756504472
$bf[0,3] = 1.045534700225387
$bf[0,4] = 1.05175480738341
$ws.Range("B12:F12").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036609496815501
$in[0,1] = 1.039049085872673
$in[0,2] = 1.040148811353839
$in[0,3] = 1.049111317484688
$in[0,4] = 1.055308674376457
$in[0,5] = 1.016944511814853
$ws.Range("I12:N12").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032422470312568
$bf[0,2] = 1.036573231929887
$bf[0,3] = 1.045581567053494
$bf[0,4] = 1.051809052668002
$ws.Range("B13:F13").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036619968135384
$in[0,1] = 1.039076438662089
$in[0,2] = 1.040173523578139
$in[0,3] = 1.049148473114462
$in[0,4] = 1.055353248559613
$in[0,5] = 1.016953708951329
$ws.Range("I13:N13").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032573157941337
$bf[0,2] = 1.036686206651518
$bf[0,3] = 1.045735161031798
$bf[0,4] = 1.051986832942865
$ws.Range("B14:F14").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036654248333551
$in[0,1] = 1.039166059083569
$in[0,2] = 1.040254486178888
$in[0,3] = 1.049270227131853
$in[0,4] = 1.055499321275833
$in[0,5] = 1.016983841510739
$ws.Range("I14:N14").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032666015242134
$bf[0,2] = 1.036755824358545
$bf[0,3] = 1.045829819309829
$bf[0,4] = 1.052096401032404
$ws.Range("B15:F15").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036675346630788
$in[0,1] = 1.039221274850425
$in[0,2] = 1.040304363126556
$in[0,3] = 1.049345252073927
$in[0,4] = 1.055589338386506
$in[0,5] = 1.017002405223026
$ws.Range("I15:N15").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033206692781392
$bf[0,2] = 1.037161184885424
$bf[0,3] = 1.046381136107931
$bf[0,4] = 1.052734619070433
$ws.Range("B16:F16").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036797798296354
$in[0,1] = 1.039542618087352
$in[0,2] = 1.040594564150905
$in[0,3] = 1.04978205506613
$in[0,4] = 1.056113533777492
$in[0,5] = 1.017110423944891
$ws.Range("I16:N16").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033546027613846
$bf[0,2] = 1.037415592436812
$bf[0,3] = 1.046727282585935
$bf[0,4] = 1.053135381871797
$ws.Range("B17:F17").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036874300666995
$in[0,1] = 1.039744155326427
$in[0,2] = 1.040776506799338
$in[0,3] = 1.050056159496528
$in[0,4] = 1.056442572910964
$in[0,5] = 1.017178154394159
$ws.Range("I17:N17").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033744021100403
$bf[0,2] = 1.037564032431735
$bf[0,3] = 1.046929298824461
$bf[0,4] = 1.053369292918481
$ws.Range("B18:F18").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036918811907514
$in[0,1] = 1.039861696408628
$in[0,2] = 1.040882597087329
$in[0,3] = 1.050216079025416
$in[0,4] = 1.056634576618494
$in[0,5] = 1.01721765061801
$ws.Range("I18:N18").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033811542786831
$bf[0,2] = 1.037614654792283
$bf[0,3] = 1.046998200577571
$bf[0,4] = 1.053449076429566
$ws.Range("B19:F19").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036933970194954
$in[0,1] = 1.039901772783151
$in[0,2] = 1.040918765371834
$in[0,3] = 1.050270614038787
$in[0,4] = 1.056700058577856
$in[0,5] = 1.017231116135072
$ws.Range("I19:N19").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033509613442642
$bf[0,2] = 1.037388291899532
$bf[0,3] = 1.046690132454151
$bf[0,4] = 1.053092367992148
$ws.Range("B20:F20").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036866104194609
$in[0,1] = 1.039722533540406
$in[0,2] = 1.040756989561463
$in[0,3] = 1.050026746637055
$in[0,4] = 1.056407261764449
$in[0,5] = 1.017170888568862
$ws.Range("I20:N20").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.03252878133844
$bf[0,2] = 1.036652936273236
$bf[0,3] = 1.04568992643958
$bf[0,4] = 1.051934474456836
$ws.Range("B21:F21").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036644158431739
$in[0,1] = 1.039139668587279
$in[0,2] = 1.040230646122002
$in[0,3] = 1.049234371882585
$in[0,4] = 1.05545630298092
$in[0,5] = 1.016974968628971
$ws.Range("I21:N21").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.031912689517502
$bf[0,2] = 1.036191035401007
$bf[0,3] = 1.045062102118062
$bf[0,4] = 1.051207847919715
$ws.Range("B22:F22").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036503613553379
$in[0,1] = 1.038773095577817
$in[0,2] = 1.039899416098033
$in[0,3] = 1.048736534676149
$in[0,4] = 1.054859133823249
$in[0,5] = 1.016851700431528
$ws.Range("I22:N22").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032239234508465
$bf[0,2] = 1.036435854966131
$bf[0,3] = 1.045394824131229
$bf[0,4] = 1.051592914441809
$ws.Range("B23:F23").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.036578213534102
$in[0,1] = 1.038967432256669
$in[0,2] = 1.040075035199212
$in[0,3] = 1.049000413041556
$in[0,4] = 1.055175634259572
$in[0,5] = 1.016917055211829
$ws.Range("I23:N23").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033526067237787
$bf[0,2] = 1.037400627687908
$bf[0,3] = 1.046706918644878
$bf[0,4] = 1.053111803640605
$ws.Range("B24:F24").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.03686980817162
$in[0,1] = 1.039732303532553
$in[0,2] = 1.040765808665644
$in[0,3] = 1.050040036921655
$in[0,4] = 1.056423217102243
$in[0,5] = 1.017174171712889
$ws.Range("I24:N24").Value = $in

$bf = New-Object 'double[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.035021918877716
$bf[0,2] = 1.038522084779374
$bf[0,3] = 1.048234031448076
$bf[0,4] = 1.054880382272506
$ws.Range("B25:F25").Value = $bf

$in = New-Object 'double[,]' 1,6
$in[0,0] = 1.037203819467953
$in[0,1] = 1.040619412584745
$in[0,2] = 1.041566083825987
$in[0,3] = 1.051247990072506
$in[0,4] = 1.057874130445436
$in[0,5] = 1.017472155167944
$ws.Range("I25:N25").Value = $in
